$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New loading_percent values for rows 2-25 (data rows 0-23), columns B,C,D,E,F,I,J,K,M
# corresponding to the "case with 380 kV done" update.

$colB = @(12.13083128559418, 12.04086693378411, 11.98964489756303, 11.9698055924706, 11.96657441124192, 11.98937312047309, 12.09899307288976, 12.34466704629653, 12.5422063228092, 12.6353848751939, 12.67111206126513, 12.66339844941499, 12.63831550426018, 12.62300806524374, 12.53618070101222, 12.48374019968291, 12.45389213952835, 12.44384112303781, 12.48929029608429, 12.64567124267105, 12.75043794322553, 12.69429906556653, 12.48678015830039, 12.27508518096457)
$colC = @(10.47083485335457, 10.37247870353628, 10.3147501515664, 10.29191967988673, 10.28817133365391, 10.31443940858476, 10.43638441402232, 10.69548322837257, 10.89637817595271, 10.98969237793403, 11.02527361381435, 11.01760016211028, 10.99261493679085, 10.97734174653709, 10.89031591838782, 10.83739916630715, 10.80714606728003, 10.79693532027974, 10.84301351227536, 10.99994730353215, 11.10392682654408, 11.0483121400207, 10.84047473624926, 10.62341790507813)
$colD = @(5.28821961223046, 5.317610910191508, 5.336914150959911, 5.345097126633478, 5.346475054978033, 5.337023225995006, 5.298093289002699, 5.23169332289188, 5.188927412982413, 5.170769980036344, 5.164080058425564, 5.165512594103984, 5.170215873780915, 5.173120960874352, 5.190140091116693, 5.200912536596503, 5.207230677030529, 5.20939088260532, 5.199753157378216, 5.168829365902186, 5.149702198813885, 5.159811800934128, 5.20027692357875, 5.248596386053745)
$colE = @(11.81046704551659, 11.82745611112445, 11.83991031635698, 11.84549423694146, 11.84645216902843, 11.83998356320584, 11.81590509864028, 11.78473225718041, 11.77159993853598, 11.76774293534281, 11.7665862341402, 11.76682184608966, 11.76764168655395, 11.76818341601086, 11.77189455655568, 11.77471306559574, 11.77653349540067, 11.7771841062829, 11.77439240759064, 11.76739263776822, 11.76458868904591, 11.76592339217849, 11.77453675397184, 11.79144835106125)
$colF = @(51.50541655708697, 51.23820060089868, 51.08074179322472, 51.01827354326347, 51.00800422687843, 51.07989240767553, 51.41192946782491, 52.11390353426017, 52.65834176497859, 52.91172146767533, 53.00844342986918, 52.98757905716525, 52.91966364308003, 52.87816262835507, 52.64189360735276, 52.49837895251304, 52.41637378431277, 52.38870256575666, 52.513600724832, 52.93959147544163, 53.22248054920482, 53.07110394484378, 52.50671739002522, 51.91878645273509)
$colI = @(35.28633954450259, 35.14514794079964, 35.06197438401576, 35.02897880520535, 35.02355454567333, 35.06152573984583, 35.23693478587316, 35.60819839588336, 35.89670185011016, 36.03115638467137, 36.08251348016773, 36.07143347448067, 36.03537282073241, 36.01334159654591, 35.8879778958191, 35.81188029302526, 35.76841461912829, 35.75375063587174, 35.81994972768111, 36.04595287351016, 36.19622601020012, 36.11579430544898, 35.81630065054787, 35.50492947291686)
$colJ = @(10.52997060057707, 10.54178900088674, 10.55010298782948, 10.55375708813996, 10.55437992443631, 10.55015119071906, 10.53382622564335, 10.51019541000403, 10.49793332263314, 10.49345969470607, 10.49192420386951, 10.49224785076447, 10.49333019235318, 10.4940138013271, 10.49824788884071, 10.5011280753793, 10.50288865529546, 10.50350262190983, 10.50081071553282, 10.49300798077412, 10.48883260940944, 10.49097660899449, 10.50095386770738, 10.51569179968744)
$colK = @(13.82379104765765, 13.77240958628189, 13.74485088046359, 13.73463357201729, 13.7329984444186, 13.7447089723594, 13.80525296986628, 13.95514050277604, 14.08350743057729, 14.14568078525626, 14.16975024953961, 14.16454336683786, 14.14765056175468, 14.13737116180953, 14.07951880970731, 14.04498386029917, 14.02547705044073, 14.01893417645637, 14.048623363614, 14.15259826411024, 14.22360759230101, 14.1854350413659, 14.04697685950987, 13.91132861611729)
$colM = @(17.03169587230134, 17.03829306940236, 17.0460973057287, 17.05022212707956, 17.05096410934267, 17.04614910958167, 17.03319206768051, 17.03752013410459, 17.05874291357842, 17.07228812916987, 17.0779737161335, 17.07672453716332, 17.07274475655916, 17.07037937090745, 17.05793573376525, 17.05129638646, 17.04784405743356, 17.04673818551209, 17.05196525205957, 17.0738986443423, 17.09147415940599, 17.08179842145176, 17.05166172162985, 17.03317275678769)

$startRow = 2

for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $colB[$i]
}

for ($i = 0; $i -lt $colC.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $colC[$i]
}

for ($i = 0; $i -lt $colD.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 4).Value = $colD[$i]
}

for ($i = 0; $i -lt $colE.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 5).Value = $colE[$i]
}

for ($i = 0; $i -lt $colF.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 6).Value = $colF[$i]
}

for ($i = 0; $i -lt $colI.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 9).Value = $colI[$i]
}

for ($i = 0; $i -lt $colJ.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 10).Value = $colJ[$i]
}

for ($i = 0; $i -lt $colK.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 11).Value = $colK[$i]
}

for ($i = 0; $i -lt $colM.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 13).Value = $colM[$i]
}
